# Scheduled runner update: refresh market-price-derived columns (H:N)
# for the rows whose underlying Universalis data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1197.2667
$ws.Range("I17").Value = 893
$ws.Range("J17").Value = 1273.3334
$ws.Range("K17").Value = 2679
$ws.Range("L17").Value = 3820.0002
$ws.Range("M17").Value = -2511
$ws.Range("N17").Value = -4156.0002

$ws.Range("H32").Value = 6336810.5
$ws.Range("I32").Value = 788
$ws.Range("J32").Value = 11616829
$ws.Range("K32").Value = 788
$ws.Range("L32").Value = 11616829
$ws.Range("M32").Value = -462
$ws.Range("N32").Value = -11617481

$ws.Range("H116").Value = 150348.92
$ws.Range("I116").Value = 208888.5
$ws.Range("K116").Value = 208888.5
$ws.Range("M116").Value = -205446.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 14019.2
$ws.Range("I63").Value = 21000
$ws.Range("J63").Value = 3548
$ws.Range("K63").Value = 21000
$ws.Range("L63").Value = 3548
$ws.Range("M63").Value = -20314
$ws.Range("N63").Value = -4920

$ws.Range("H66").Value = 14019.2
$ws.Range("I66").Value = 21000
$ws.Range("J66").Value = 3548
$ws.Range("K66").Value = 105000
$ws.Range("L66").Value = 17740
$ws.Range("M66").Value = -101568
$ws.Range("N66").Value = -24604

$ws.Range("H132").Value = 3869
$ws.Range("I132").Value = 3883.087
$ws.Range("J132").Value = 3845.8572
$ws.Range("K132").Value = 11649.261
$ws.Range("L132").Value = 11537.5716
$ws.Range("M132").Value = -9119.261
$ws.Range("N132").Value = -16597.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 28724
$ws.Range("J40").Value = 28724
$ws.Range("L40").Value = 28724
$ws.Range("N40").Value = -29254

$ws.Range("H86").Value = 6152.364
$ws.Range("I86").Value = 4528.25
$ws.Range("J86").Value = 10483.333
$ws.Range("K86").Value = 4528.25
$ws.Range("L86").Value = 10483.333
$ws.Range("M86").Value = -3405.25
$ws.Range("N86").Value = -12729.333

$ws.Range("H87").Value = 19857.143
$ws.Range("J87").Value = 19857.143
$ws.Range("L87").Value = 19857.143
$ws.Range("N87").Value = -22353.143

$ws.Range("H89").Value = 6152.364
$ws.Range("I89").Value = 4528.25
$ws.Range("J89").Value = 10483.333
$ws.Range("K89").Value = 22641.25
$ws.Range("L89").Value = 52416.665
$ws.Range("M89").Value = -17025.25
$ws.Range("N89").Value = -63648.665

$ws.Range("H90").Value = 19857.143
$ws.Range("J90").Value = 19857.143
$ws.Range("L90").Value = 59571.429
$ws.Range("N90").Value = -72051.429

$ws.Range("H96").Value = 15065.6
$ws.Range("I96").Value = 5109.3335
$ws.Range("K96").Value = 5109.3335
$ws.Range("M96").Value = -2363.3335

$ws.Range("H134").Value = 34943.516
$ws.Range("J134").Value = 6883.3335
$ws.Range("L134").Value = 20650.0005
$ws.Range("N134").Value = -25720.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 46550
$ws.Range("J92").Value = 46550
$ws.Range("L92").Value = 46550
$ws.Range("N92").Value = -51542

$ws.Range("H122").Value = 1163.238
$ws.Range("I122").Value = 825.2308
$ws.Range("J122").Value = 1712.5
$ws.Range("K122").Value = 2475.6924
$ws.Range("L122").Value = 5137.5
$ws.Range("M122").Value = -25.69239999999991
$ws.Range("N122").Value = -10037.5

$ws.Range("H132").Value = 2811.1428
$ws.Range("I132").Value = 1522.2
$ws.Range("J132").Value = 3982.9092
$ws.Range("K132").Value = 4566.6
$ws.Range("L132").Value = 11948.7276
$ws.Range("M132").Value = -2036.6
$ws.Range("N132").Value = -17008.7276

$ws.Range("H134").Value = 2292.7407
$ws.Range("I134").Value = 1371
$ws.Range("J134").Value = 2834.9412
$ws.Range("K134").Value = 4113
$ws.Range("L134").Value = 8504.8236
$ws.Range("M134").Value = -1578
$ws.Range("N134").Value = -13574.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1981.8
$ws.Range("I109").Value = 972.7
$ws.Range("J109").Value = 4000
$ws.Range("K109").Value = 2918.1
$ws.Range("L109").Value = 12000
$ws.Range("M109").Value = -1878.1
$ws.Range("N109").Value = -14080

$ws.Range("H131").Value = 890.25354
$ws.Range("J131").Value = 906.91174
$ws.Range("L131").Value = 2720.73522
$ws.Range("N131").Value = -12800.73522

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 3294
$ws.Range("I99").Value = 1490
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 1490
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = 756
$ws.Range("N99").Value = -10492

$ws.Range("H132").Value = 5196.6665
$ws.Range("I132").Value = 6273.7144
$ws.Range("K132").Value = 18821.1432
$ws.Range("M132").Value = -16291.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1930.7778
$ws.Range("I93").Value = 2045.8182
$ws.Range("J93").Value = 1750
$ws.Range("K93").Value = 2045.8182
$ws.Range("L93").Value = 1750
$ws.Range("M93").Value = -797.8181999999999
$ws.Range("N93").Value = -4246

$ws.Range("H136").Value = 5230.15
$ws.Range("I136").Value = 2947.0908
$ws.Range("J136").Value = 8020.5557
$ws.Range("K136").Value = 8841.2724
$ws.Range("L136").Value = 24061.6671
$ws.Range("M136").Value = -6291.2724
$ws.Range("N136").Value = -29161.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 26433.6
$ws.Range("J16").Value = 26433.6
$ws.Range("L16").Value = 26433.6
$ws.Range("N16").Value = -27017.6

$ws.Range("H75").Value = 26999.5
$ws.Range("I75").Value = 14999
$ws.Range("J75").Value = 39000
$ws.Range("K75").Value = 14999
$ws.Range("L75").Value = 39000
$ws.Range("M75").Value = -14063
$ws.Range("N75").Value = -40872

$ws.Range("H78").Value = 26999.5
$ws.Range("I78").Value = 14999
$ws.Range("J78").Value = 39000
$ws.Range("K78").Value = 44997
$ws.Range("L78").Value = 117000
$ws.Range("M78").Value = -40317
$ws.Range("N78").Value = -126360

$ws.Range("H93").Value = 13133.333
$ws.Range("I93").Value = 9000
$ws.Range("J93").Value = 13960
$ws.Range("K93").Value = 9000
$ws.Range("L93").Value = 13960
$ws.Range("M93").Value = -6504
$ws.Range("N93").Value = -18952

$ws.Range("H94").Value = 16168.333
$ws.Range("I94").Value = 9800
$ws.Range("J94").Value = 16747.273
$ws.Range("K94").Value = 9800
$ws.Range("L94").Value = 16747.273
$ws.Range("M94").Value = -8899
$ws.Range("N94").Value = -18549.273

$ws.Range("H110").Value = 25822
$ws.Range("J110").Value = 25822
$ws.Range("L110").Value = 25822
$ws.Range("N110").Value = -34002

$ws.Range("H116").Value = 27493.334
$ws.Range("J116").Value = 27493.334
$ws.Range("L116").Value = 27493.334
$ws.Range("N116").Value = -36671.334

$ws.Range("H132").Value = 27074.977
$ws.Range("I132").Value = 113182.445
$ws.Range("J132").Value = 2857.25
$ws.Range("K132").Value = 339547.335
$ws.Range("L132").Value = 8571.75
$ws.Range("M132").Value = -337017.335
$ws.Range("N132").Value = -13631.75
